{"js": "// Applies the \"Change placeTile1 test cases and updated doc\" edit.\n//\n// Four textual edits inside the test-case list, each implemented as a\n// find -> replace (the replace collapses the touched runs into one,\n// which mirrors how the target document merges the split runs that\n// used to make up this text):\n//\n//   1. \"placeTile\" / \"1\" / \":\"                                  -> \"placeTile1:\"\n//   2. \"3. Improper \" / \"placing of tile\" / ... / \"C10\"         -> \"3. Improper placing of tile (on top of another tile)-> place O1 at C10\"\n//   3. \"...fails -> place R5 at M5\"                             -> \"...fails -> place R5 at L6\"\n//   4. \"...fails ->\" / \"place O1 at L6 ...\" / \" place R5 at F10\" / \" ( checking ...)\"\n//                                                                -> \"...fails ->place O1 at L6 (checking next to one tile)\"\n\nasync function findUnique(context, text) {\n  const results = context.document.body.search(text, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  return results.items[0];\n}\n\n// 1) \"placeTile\" + \"1\" + \":\" -> single run \"placeTile1:\"\nconst r1 = await findUnique(context, \"placeTile1:\");\nr1.insertText(\"placeTile1:\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Merge the \"3. Improper ... C10\" runs, text itself is unchanged.\nconst r2 = await findUnique(\n  context,\n  \"3. Improper placing of tile (on top of another tile)-> place O1 at C10\"\n);\nr2.insertText(\n  \"3. Improper placing of tile (on top of another tile)-> place O1 at C10\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 3) \"place R5 at M5\" -> \"place R5 at L6\"\nconst r3 = await findUnique(context, \"M5\");\nr3.insertText(\"L6\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 4) Drop the trailing \"place R5 at F10 ( checking next to a tile in a segment)\"\n//    clause and merge what remains into a single run.\nconst r4 = await findUnique(\n  context,\n  \"Placing a tile next to an (adjacent) duplicate fails ->place O1 at L6 (checking next to one tile), place R5 at F10 ( checking next to a tile in a segment)\"\n);\nr4.insertText(\n  \"Placing a tile next to an (adjacent) duplicate fails ->place O1 at L6 (checking next to one tile)\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# Applies the \"Change placeTile1 test cases and updated doc\" edit.\n#\n# Four textual edits inside the test-case list, each implemented as a\n# Find/Replace (wdReplaceAll = 2), which collapses the touched runs\n# into one - matching how the target document merges the runs that\n# used to make up this text:\n#\n#   1. \"placeTile\" / \"1\" / \":\"                                  -> \"placeTile1:\"\n#   2. \"3. Improper \" / \"placing of tile\" / ... / \"C10\"         -> \"3. Improper placing of tile (on top of another tile)-> place O1 at C10\"\n#   3. \"...fails -> place R5 at M5\"                             -> \"...fails -> place R5 at L6\"\n#   4. \"...fails ->\" / \"place O1 at L6 ...\" / \" place R5 at F10\" / \" ( checking ...)\"\n#                                                                -> \"...fails ->place O1 at L6 (checking next to one tile)\"\n\n$d = $word.ActiveDocument\n\nfunction ReplaceText($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# 1) \"placeTile\" + \"1\" + \":\" -> single run \"placeTile1:\"\nReplaceText \"placeTile1:\" \"placeTile1:\"\n\n# 2) Merge the \"3. Improper ... C10\" runs, text itself is unchanged.\nReplaceText \"3. Improper placing of tile (on top of another tile)-> place O1 at C10\" \"3. Improper placing of tile (on top of another tile)-> place O1 at C10\"\n\n# 3) \"place R5 at M5\" -> \"place R5 at L6\"\nReplaceText \"M5\" \"L6\"\n\n# 4) Drop the trailing \"place R5 at F10 ( checking next to a tile in a segment)\"\n#    clause and merge what remains into a single run.\nReplaceText \"Placing a tile next to an (adjacent) duplicate fails ->place O1 at L6 (checking next to one tile), place R5 at F10 ( checking next to a tile in a segment)\" \"Placing a tile next to an (adjacent) duplicate fails ->place O1 at L6 (checking next to one tile)\"\n"}
